$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4:B73 with the corrected naive-forecaster values
$ws.Cells.Item(4, 2).Value = -0.3499999999999943
$ws.Cells.Item(5, 2).Value = 0.25
$ws.Cells.Item(6, 2).Value = -0.8999999999999915
$ws.Cells.Item(7, 2).Value = -0.6500000000000057
$ws.Cells.Item(8, 2).Value = -0.7000000000000171
$ws.Cells.Item(9, 2).Value = 0.2000000000000028
$ws.Cells.Item(10, 2).Value = 0.4999999999999858
$ws.Cells.Item(11, 2).Value = 0.4200000000000017
$ws.Cells.Item(12, 2).Value = 1.099999999999994
$ws.Cells.Item(13, 2).Value = 0.4999999999999858
$ws.Cells.Item(14, 2).Value = 0.8
$ws.Cells.Item(15, 2).Value = 0.3
$ws.Cells.Item(16, 2).Value = 0.4
$ws.Cells.Item(17, 2).Value = 0.5
$ws.Cells.Item(18, 2).Value = -0.2
$ws.Cells.Item(19, 2).Value = 0.1
$ws.Cells.Item(20, 2).Value = 0.1
$ws.Cells.Item(21, 2).Value = 0.1
$ws.Cells.Item(22, 2).Value = -0.3
$ws.Cells.Item(23, 2).Value = 0.2
$ws.Cells.Item(24, 2).Value = 0.8999999999999915
$ws.Cells.Item(25, 2).Value = 0.4200000000000017
$ws.Cells.Item(26, 2).Value = 0.3
$ws.Cells.Item(27, 2).Value = 0.5
$ws.Cells.Item(28, 2).Value = 0.2999999999999829
$ws.Cells.Item(29, 2).Value = 0.4000000000000057
$ws.Cells.Item(30, 2).Value = 0.2000000000000028
$ws.Cells.Item(31, 2).Value = 0.4200000000000017
$ws.Cells.Item(32, 2).Value = 0.6200000000000045
$ws.Cells.Item(33, 2).Value = 0.4200000000000017
$ws.Cells.Item(34, 2).Value = 0.3000000000000114
$ws.Cells.Item(35, 2).Value = 0.5400000000000063
$ws.Cells.Item(36, 2).Value = 0.3400000000000034
$ws.Cells.Item(37, 2).Value = 0.4399999999999977
$ws.Cells.Item(38, 2).Value = 0.4999999999999858
$ws.Cells.Item(39, 2).Value = 0.4999999999999858
$ws.Cells.Item(40, 2).Value = 0.6999999999999886
$ws.Cells.Item(41, 2).Value = 0.5999999999999943
$ws.Cells.Item(42, 2).Value = 0.5999999999999943
$ws.Cells.Item(43, 2).Value = 0.6999999999999886
$ws.Cells.Item(44, 2).Value = 0.2999999999999829
$ws.Cells.Item(45, 2).Value = 0.4999999999999716
$ws.Cells.Item(46, 2).Value = 0.2
$ws.Cells.Item(47, 2).Value = 0.09999999999999432
$ws.Cells.Item(48, 2).Value = -0.1
$ws.Cells.Item(49, 2).Value = -0.09999999999999432
$ws.Cells.Item(50, 2).Value = 0.08000000000004093
$ws.Cells.Item(51, 2).Value = 0.09999999999999432
$ws.Cells.Item(52, 2).Value = -11.9
$ws.Cells.Item(53, 2).Value = 6.640000000000001
$ws.Cells.Item(54, 2).Value = -0.4000000000000057
$ws.Cells.Item(55, 2).Value = -0.7094799999999992
$ws.Cells.Item(56, 2).Value = 1.310000000000016
$ws.Cells.Item(57, 2).Value = 1.52000000000001
$ws.Cells.Item(58, 2).Value = -0.539999999999992
$ws.Cells.Item(59, 2).Value = 0.4652855479103435
$ws.Cells.Item(60, 2).Value = 0.38
$ws.Cells.Item(61, 2).Value = 0.04
$ws.Cells.Item(62, 2).Value = -0.29
$ws.Cells.Item(63, 2).Value = -0.2078779574152918
$ws.Cells.Item(64, 2).Value = 0.1206478331785803
$ws.Cells.Item(65, 2).Value = -0.18
$ws.Cells.Item(66, 2).Value = 0.044
$ws.Cells.Item(67, 2).Value = -0.08251004046350374
$ws.Cells.Item(68, 2).Value = 0.2582525219575302
$ws.Cells.Item(69, 2).Value = -0.04717552522494373
$ws.Cells.Item(70, 2).Value = 0.2142297805489477
$ws.Cells.Item(71, 2).Value = 0.2394371574146135
$ws.Cells.Item(72, 2).Value = 0.04717883418304325
$ws.Cells.Item(73, 2).Value = 0.0959495356205764

# Remove the now-stale forecast rows 74-82 (data series shortened)
$ws.Range("A74:B82").EntireRow.Delete()
